# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 236-241 (columns B, F-AC). Column A (id), C, D, E stay the same.
$rows = @(
    @{ Row = 236; B = 6012013; F = "Valerenga"; G = "Tromso"; H = 1; I = 1; J = "D";
       K = 2.5; L = 3.5; M = 2.625; N = 2.05; O = 3.6; P = 3.4; Q = -0.25; R = 1.8;
       S = 2.05; T = 2.5; U = 2; V = 1.85; W = -1; X = 2.6; Y = -1; Z = -0.5;
       AA = 0.5249999999999999; AB = -1; AC = 0.8500000000000001 },
    @{ Row = 237; B = 6012008; F = "Stromsgodset"; G = "SK Brann"; H = 3; I = 0; J = "H";
       K = 4.5; L = 4.2; M = 1.65; N = 4.2; O = 3.8; P = 1.75; Q = 0.75; R = 1.86;
       S = 2.04; T = 2.75; U = 1.85; V = 2; W = 3.2; X = -1; Y = -1; Z = 0.8600000000000001;
       AA = -1; AB = 0.425; AC = -0.5 },
    @{ Row = 238; B = 6011940; F = "BodoGlimt"; G = "Sarpsborg"; H = 2; I = 0; J = "H";
       K = 1.45; L = 5; M = 5.75; N = 1.45; O = 5; P = 6; Q = -1.25; R = 1.89;
       S = 2.01; T = 4; U = 1.85; V = 2; W = 0.45; X = -1; Y = -1; Z = 0.8899999999999999;
       AA = -1; AB = -1; AC = 1 },
    @{ Row = 239; B = 6011535; F = "Odd BK"; G = "Aalesund"; H = 4; I = 1; J = "H";
       K = 1.6; L = 4.2; M = 5; N = 1.5; O = 4.75; P = 5.75; Q = -1.25; R = 2.05;
       S = 1.8; T = 3.25; U = 2.025; V = 1.825; W = 0.5; X = -1; Y = -1; Z = 1.05;
       AA = -1; AB = 1.025; AC = -1 },
    @{ Row = 240; B = 6011534; F = "Molde"; G = "HamKam"; H = 1; I = 1; J = "D";
       K = 1.25; L = 6.5; M = 10; N = 1.222; O = 7; P = 11; Q = -2; R = 1.95;
       S = 1.9; T = 3.5; U = 1.8; V = 2; W = -1; X = 6; Y = -1; Z = -1;
       AA = 0.8999999999999999; AB = -1; AC = 1 },
    @{ Row = 241; B = 6390445; F = "Haugesund"; G = "Stabaek"; H = 3; I = 0; J = "H";
       K = 2.25; L = 3.4; M = 3.1; N = 1.8; O = 3.75; P = 4.5; Q = -0.75; R = 1.975;
       S = 1.875; T = 2.5; U = 1.9; V = 1.95; W = 0.8; X = -1; Y = -1; Z = 0.9750000000000001;
       AA = -1; AB = 0.8999999999999999; AC = -1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = $r.B
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W
    $ws.Range("X$row").Value = $r.X
    $ws.Range("Y$row").Value = $r.Y
    $ws.Range("Z$row").Value = $r.Z
    $ws.Range("AA$row").Value = $r.AA
    $ws.Range("AB$row").Value = $r.AB
    $ws.Range("AC$row").Value = $r.AC
}
